# Update bus-arrival schedule data (commit: "Arribos 141 actualizados - 31").
# Rows 5+ on sheet "TODOS" / "COMBINADAS" and rows 2+ on sheet "215" are
# rewritten in place with refreshed ETA/BANDERA/MIN/ESTADO values, and the
# schedule grows from 49 to 68 rows (TODOS / COMBINADAS) and 9 to 12 rows (215).

$wb = $excel.ActiveWorkbook

# --- Sheet "TODOS" ---
$wsTodos = $wb.Worksheets.Item("TODOS")
$wsTodosData = @(
    @(5, "19:36", "16_P MOR-SANTA ANA", 0, "🚌"),
    @(6, "19:36", "16_SANTA ANA", 0, "🚌"),
    @(7, "19:36", "16_SANTA ANA", 0, "🚌"),
    @(8, "19:36", "23_HERNANDEZ", 0, "🚌"),
    @(9, "18:31", "15_ABASTO", 1, "🚌"),
    @(10, "19:37", "15_ABASTO", 1, "🚌"),
    @(11, "18:35", "23_HERNANDEZ", 5, "🚌"),
    @(12, "19:41", "14_ABASTO", 5, "🚌"),
    @(13, "18:40", "14_ABASTO", 10, "📅"),
    @(14, "18:40", "15_ABASTO", 10, "🚌"),
    @(15, "19:51", "16_P MOR-SANTA ANA", 15, "🚌"),
    @(16, "19:00", "16_SANTA ANA", 17, "🚌"),
    @(17, "18:48", "14X44_ABASTO", 18, "🚌"),
    @(18, "19:04", "23_HERNANDEZ", 21, "🚌"),
    @(19, "18:52", "215A_LA PLATA", 22, "🚌"),
    @(20, "20:01", "14_ABASTO", 25, "🚌"),
    @(21, "18:56", "10_OLMOS", 26, "🚌"),
    @(22, "19:10", "14_ABASTO", 27, "🚌"),
    @(23, "18:58", "215A_EL PATO", 28, "📅"),
    @(24, "19:12", "215B_LP-P MOR-1 Y 57", 29, "🚌"),
    @(25, "19:30", "16_SANTA ANA", 30, "🚌"),
    @(26, "19:01", "16_SANTA ANA", 31, "🚌"),
    @(27, "19:16", "15_ABASTO", 33, "🚌"),
    @(28, "19:16", "27_EL RETIRO", 33, "📅"),
    @(29, "19:04", "11_ETCHEVERRY", 34, "🚌"),
    @(30, "19:05", "23_HERNANDEZ", 35, "🚌"),
    @(31, "20:11", "10_OLMOS", 35, "🚌"),
    @(32, "20:11", "16_P MOR-167 Y 521", 35, "📅"),
    @(33, "19:20", "14_ABASTO", 37, "📅"),
    @(34, "19:20", "16_SANTA ANA", 37, "🚌"),
    @(35, "20:13", "23_HERNANDEZ", 37, "🚌"),
    @(36, "19:10", "16_P MOR-SANTA ANA", 40, "🚌"),
    @(37, "19:10", "215B_LP-P MOR-1 Y 57", 40, "🚌"),
    @(38, "19:12", "10_OLMOS", 42, "🚌"),
    @(39, "19:28", "15_ABASTO", 45, "🚌"),
    @(40, "19:17", "27_EL RETIRO", 47, "🚌"),
    @(41, "20:24", "215A_EL PATO", 48, "📅"),
    @(42, "19:21", "16_SANTA ANA", 51, "🚌"),
    @(43, "19:21", "26_HERNANDEZ", 51, "🚌"),
    @(44, "19:34", "23_HERNANDEZ", 51, "🚌"),
    @(45, "19:30", "225_GOMEZ", 60, "📅"),
    @(46, "20:44", "11_ETCHEVERRY", 68, "🚌"),
    @(47, "20:09", "15_ABASTO", 69, "🚌"),
    @(48, "19:40", "14_ABASTO", 70, "🚌"),
    @(49, "19:40", "215C_EL PATO", 70, "🚌"),
    @(50, "20:10", "10_OLMOS", 70, "🚌"),
    @(51, "19:50", "11X44_ETCHEVERRY", 80, "🚌"),
    @(52, "19:50", "16_P MOR-SANTA ANA", 80, "🚌"),
    @(53, "20:56", "10_OLMOS", 80, "🚌"),
    @(54, "19:51", "81_EL PELIGRO", 81, "🚌"),
    @(55, "20:57", "27_EL RETIRO", 81, "🚌"),
    @(56, "19:54", "215C_LA PLATA", 84, "🚌"),
    @(57, "21:04", "84_COLONIA URQUIZA-ESC 49", 88, "🚌"),
    @(58, "19:59", "17_ROMERO", 89, "📅"),
    @(59, "21:08", "215B_EL PATO", 92, "🚌"),
    @(60, "20:10", "16_P MOR-167 Y 521", 100, "🚌"),
    @(61, "21:21", "26_HERNANDEZ", 105, "📅"),
    @(62, "20:31", "225_GOMEZ", 108, "📅"),
    @(63, "20:21", "26_HERNANDEZ", 111, "🚌"),
    @(64, "20:22", "11_ETCHEVERRY", 112, "🚌"),
    @(65, "20:23", "215A_EL PATO", 113, "🚌"),
    @(66, "21:29", "215C_LA PLATA", 113, "🚌"),
    @(67, "20:39", "215A_LA PLATA", 116, "📅"),
    @(68, "20:56", "27_EL RETIRO", 116, "🚌")
)
foreach ($row in $wsTodosData) {
    $r = $row[0]
    $wsTodos.Cells.Item($r,1).Value = $row[1]
    $wsTodos.Cells.Item($r,2).Value = $row[2]
    $wsTodos.Cells.Item($r,3).Value = $row[3]
    $wsTodos.Cells.Item($r,4).Value = $row[4]
}

# --- Sheet "215" ---
$ws215 = $wb.Worksheets.Item("215")
$ws215Data = @(
    @(2, "18:52", "215A_LA PLATA", 22, "🚌"),
    @(3, "18:58", "215A_EL PATO", 28, "📅"),
    @(4, "19:12", "215B_LP-P MOR-1 Y 57", 29, "🚌"),
    @(5, "19:10", "215B_LP-P MOR-1 Y 57", 40, "🚌"),
    @(6, "20:24", "215A_EL PATO", 48, "📅"),
    @(7, "19:40", "215C_EL PATO", 70, "🚌"),
    @(8, "19:54", "215C_LA PLATA", 84, "🚌"),
    @(9, "21:08", "215B_EL PATO", 92, "🚌"),
    @(10, "20:23", "215A_EL PATO", 113, "🚌"),
    @(11, "21:29", "215C_LA PLATA", 113, "🚌"),
    @(12, "20:39", "215A_LA PLATA", 116, "📅")
)
foreach ($row in $ws215Data) {
    $r = $row[0]
    $ws215.Cells.Item($r,1).Value = $row[1]
    $ws215.Cells.Item($r,2).Value = $row[2]
    $ws215.Cells.Item($r,3).Value = $row[3]
    $ws215.Cells.Item($r,4).Value = $row[4]
}

# --- Sheet "COMBINADAS" ---
$wsComb = $wb.Worksheets.Item("COMBINADAS")
$wsCombData = @(
    @(5, "19:36", "16_P MOR-SANTA ANA", 0, "🚌"),
    @(6, "19:36", "16_SANTA ANA", 0, "🚌"),
    @(7, "19:36", "16_SANTA ANA", 0, "🚌"),
    @(8, "19:36", "23_HERNANDEZ", 0, "🚌"),
    @(9, "18:31", "15_ABASTO", 1, "🚌"),
    @(10, "19:37", "15_ABASTO", 1, "🚌"),
    @(11, "18:35", "23_HERNANDEZ", 5, "🚌"),
    @(12, "19:41", "14_ABASTO", 5, "🚌"),
    @(13, "18:40", "14_ABASTO", 10, "📅"),
    @(14, "18:40", "15_ABASTO", 10, "🚌"),
    @(15, "19:51", "16_P MOR-SANTA ANA", 15, "🚌"),
    @(16, "19:00", "16_SANTA ANA", 17, "🚌"),
    @(17, "18:48", "14X44_ABASTO", 18, "🚌"),
    @(18, "19:04", "23_HERNANDEZ", 21, "🚌"),
    @(19, "18:52", "215A_LA PLATA", 22, "🚌"),
    @(20, "20:01", "14_ABASTO", 25, "🚌"),
    @(21, "18:56", "10_OLMOS", 26, "🚌"),
    @(22, "19:10", "14_ABASTO", 27, "🚌"),
    @(23, "18:58", "215A_EL PATO", 28, "📅"),
    @(24, "19:12", "215B_LP-P MOR-1 Y 57", 29, "🚌"),
    @(25, "19:30", "16_SANTA ANA", 30, "🚌"),
    @(26, "19:01", "16_SANTA ANA", 31, "🚌"),
    @(27, "19:16", "15_ABASTO", 33, "🚌"),
    @(28, "19:16", "27_EL RETIRO", 33, "📅"),
    @(29, "19:04", "11_ETCHEVERRY", 34, "🚌"),
    @(30, "19:05", "23_HERNANDEZ", 35, "🚌"),
    @(31, "20:11", "10_OLMOS", 35, "🚌"),
    @(32, "20:11", "16_P MOR-167 Y 521", 35, "📅"),
    @(33, "19:20", "14_ABASTO", 37, "📅"),
    @(34, "19:20", "16_SANTA ANA", 37, "🚌"),
    @(35, "20:13", "23_HERNANDEZ", 37, "🚌"),
    @(36, "19:10", "16_P MOR-SANTA ANA", 40, "🚌"),
    @(37, "19:10", "215B_LP-P MOR-1 Y 57", 40, "🚌"),
    @(38, "19:12", "10_OLMOS", 42, "🚌"),
    @(39, "19:28", "15_ABASTO", 45, "🚌"),
    @(40, "19:17", "27_EL RETIRO", 47, "🚌"),
    @(41, "20:24", "215A_EL PATO", 48, "📅"),
    @(42, "19:21", "16_SANTA ANA", 51, "🚌"),
    @(43, "19:21", "26_HERNANDEZ", 51, "🚌"),
    @(44, "19:34", "23_HERNANDEZ", 51, "🚌"),
    @(45, "19:30", "225_GOMEZ", 60, "📅"),
    @(46, "20:44", "11_ETCHEVERRY", 68, "🚌"),
    @(47, "20:09", "15_ABASTO", 69, "🚌"),
    @(48, "19:40", "14_ABASTO", 70, "🚌"),
    @(49, "19:40", "215C_EL PATO", 70, "🚌"),
    @(50, "20:10", "10_OLMOS", 70, "🚌"),
    @(51, "19:50", "11X44_ETCHEVERRY", 80, "🚌"),
    @(52, "19:50", "16_P MOR-SANTA ANA", 80, "🚌"),
    @(53, "20:56", "10_OLMOS", 80, "🚌"),
    @(54, "19:51", "81_EL PELIGRO", 81, "🚌"),
    @(55, "20:57", "27_EL RETIRO", 81, "🚌"),
    @(56, "19:54", "215C_LA PLATA", 84, "🚌"),
    @(57, "21:04", "84_COLONIA URQUIZA-ESC 49", 88, "🚌"),
    @(58, "19:59", "17_ROMERO", 89, "📅"),
    @(59, "21:08", "215B_EL PATO", 92, "🚌"),
    @(60, "20:10", "16_P MOR-167 Y 521", 100, "🚌"),
    @(61, "21:21", "26_HERNANDEZ", 105, "📅"),
    @(62, "20:31", "225_GOMEZ", 108, "📅"),
    @(63, "20:21", "26_HERNANDEZ", 111, "🚌"),
    @(64, "20:22", "11_ETCHEVERRY", 112, "🚌"),
    @(65, "20:23", "215A_EL PATO", 113, "🚌"),
    @(66, "21:29", "215C_LA PLATA", 113, "🚌"),
    @(67, "20:39", "215A_LA PLATA", 116, "📅"),
    @(68, "20:56", "27_EL RETIRO", 116, "🚌")
)
foreach ($row in $wsCombData) {
    $r = $row[0]
    $wsComb.Cells.Item($r,1).Value = $row[1]
    $wsComb.Cells.Item($r,2).Value = $row[2]
    $wsComb.Cells.Item($r,3).Value = $row[3]
    $wsComb.Cells.Item($r,4).Value = $row[4]
}

Write-Output "Schedules updated: TODOS/COMBINADAS -> A1:D68, 215 -> A1:D12"
